# WorkPlan.xlsx update: "Work Plan completed? till christmas"
#
# Inserts two new rows into the Development block (new rows 18 & 19),
# which pushes the Misc / Other Commitments blocks down by two rows,
# fills the newly freed and previously-blank rows with the new GA
# planning tasks, relocates the existing review comment, appends two
# blank rows at the bottom of the sheet and updates the remembered
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert two rows at the top of the "Development" block.
#    Excel shifts every row from 18 downward by two (Misc header 18->20,
#    Install software 19->21, ... Easter Holiday 26->28) and keeps the
#    two already-blank trailing rows, so the sheet grows from 37 to 39
#    rows with two fresh blank rows appended at the bottom automatically.
# ---------------------------------------------------------------------
$ws.Rows("18:19").Insert()

# ---------------------------------------------------------------------
# 2. New task labels. Set these in the same order the original author
#    typed them so the shared-string table is rebuilt with the expected
#    ordering.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Research Representation"
$ws.Range("A15").Value = "Representation"
$ws.Range("A18").Value = "Basic Methods for GA"
$ws.Range("A16").Value = "Data import"
$ws.Range("A17").Value = "Validation Method"
$ws.Range("B18").Value = "Needs a slot"

# ---------------------------------------------------------------------
# 3. Formatting: the new label cells reuse the ordinary row style (A4 is
#    a plain, unbolded task-label cell); the new Gantt "bar" cells reuse
#    the existing fill styles from elsewhere in the sheet. Copy/PasteSpecial
#    (formats only) so each destination lands on the very same style index
#    as its donor cell rather than cloning a brand new one. PasteSpecial
#    only honours the first area of a multi-area destination range, so
#    paste into each target cell individually instead.
# ---------------------------------------------------------------------
$ws.Range("A4").Copy()
foreach ($addr in @("A18", "A19")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("F9").Copy()
foreach ($addr in @("H11")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("G14").Copy()
foreach ($addr in @("H15", "I15", "I16", "J16", "J17", "K17")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$ws.Range("C5").Copy()
foreach ($addr in @("H16", "H17", "I17", "J18", "K18")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. The review comment belongs to the "Autumn Exams" row, which has
#    just shifted from row 25 to row 27 - recreate it there (row-insert
#    does not relocate the legacy comment anchor automatically) and
#    drop the now-stale one.
# ---------------------------------------------------------------------
$oldComment = $ws.Range("A25").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()
$ws.Range("A27").AddComment($commentText)

# ---------------------------------------------------------------------
# 5. Restore the remembered selection on the frozen bottom-right pane,
#    which the author left on B19 after the edit.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B19").Select()
